$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G holds "K" values (formerly "Strike#"). Update rows 2-5 per regen.
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 4
$ws.Range("G5").Value = 0
